# Added test case for Successfull Payment NoCFCorp
#
# 1) Update existing "PayNowNoCFPC" sheet row-2 result data (A2/B2/K2).
# 2) Add "PayNowNoCFPS" sheet (copy of PayNowNoCFPC's layout) with its own
#    pass/fail row-2 + row-3 result data.
# 3) Add "PayNowNoCFCorp" sheet (same layout) with its own result data, and
#    make it the active sheet/tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) PayNowNoCFPC (sheet1) - refresh the latest run's result row
# ---------------------------------------------------------------------
$ws1.Range("A2").Value = "Pass"
$ws1.Range("B2").Value = "Tue Oct 22 12:40:38 IST 2024"
$ws1.Range("K2").Value = "2"

$ws1.Range("K5").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) PayNowNoCFPS (new sheet, placed right after PayNowNoCFPC)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "PayNowNoCFPS"

$ws1.Range("A1:N2").Copy()
$ws2.Range("A1:N2").PasteSpecial(-4122) | Out-Null
$ws1.Range("A1:N2").Copy()
$ws2.Range("A1:N2").PasteSpecial(-4163) | Out-Null
$ws2.Range("A1:N1").RowHeight = 29

$ws1.Range("C2:N2").Copy()
$ws2.Range("C3:N3").PasteSpecial(-4122) | Out-Null
$ws1.Range("C2").Copy()
$ws2.Range("A3:B3").PasteSpecial(-4122) | Out-Null

$ws2.Range("A2").Value = "Pass"
$ws2.Range("B2").Value = "Tue Oct 22 12:57:02 IST 2024"
$ws2.Range("D2").Value = "Y"
$ws2.Range("E2").Value = "8"
$ws2.Range("F2").Value = "880"
$ws2.Range("G2").Value = "2.5"
$ws2.Range("H2").Value = "10.50"
$ws2.Range("I2").Value = "1"
$ws2.Range("J2").Value = "1"
$ws2.Range("K2").Value = "1"
$ws2.Range("M2").Value = "3"
$ws2.Range("N2").Value = "3"

$ws2.Range("D3").Value = "N"
$ws2.Range("E3").Value = "8"
$ws2.Range("F3").Value = "880"
$ws2.Range("G3").Value = "2.5"
$ws2.Range("H3").Value = "10.50"
$ws2.Range("I3").Value = "1"
$ws2.Range("J3").Value = "1"
$ws2.Range("K3").Value = "1"
$ws2.Range("M3").Value = "3"
$ws2.Range("N3").Value = "3"

$ws2.Range("K6").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) PayNowNoCFCorp (new sheet, placed right after PayNowNoCFPS)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PayNowNoCFCorp"

$ws1.Range("A1:N2").Copy()
$ws3.Range("A1:N2").PasteSpecial(-4122) | Out-Null
$ws1.Range("A1:N2").Copy()
$ws3.Range("A1:N2").PasteSpecial(-4163) | Out-Null
$ws3.Range("A1:N1").RowHeight = 29

$ws1.Range("C2:N2").Copy()
$ws3.Range("C3:N3").PasteSpecial(-4122) | Out-Null
$ws1.Range("C2").Copy()
$ws3.Range("A3:B3").PasteSpecial(-4122) | Out-Null

$ws3.Range("A2").Value = "Pass"
$ws3.Range("B2").Value = "Wed Oct 30 15:43:35 IST 2024"
$ws3.Range("D2").Value = "Y"
$ws3.Range("E2").Value = "8"
$ws3.Range("F2").Value = "880"
$ws3.Range("G2").Value = "2.5"
$ws3.Range("H2").Value = "10.50"
$ws3.Range("I2").Value = "1"
$ws3.Range("J2").Value = "3"
$ws3.Range("K2").Value = "3"
$ws3.Range("M2").Value = "3"
$ws3.Range("N2").Value = "3"

$ws3.Range("D3").Value = "N"
$ws3.Range("E3").Value = "8"
$ws3.Range("F3").Value = "880"
$ws3.Range("G3").Value = "2.5"
$ws3.Range("H3").Value = "10.50"
$ws3.Range("I3").Value = "1"
$ws3.Range("J3").Value = "3"
$ws3.Range("K3").Value = "3"
$ws3.Range("M3").Value = "3"
$ws3.Range("N3").Value = "3"

$ws3.Range("J6").Select() | Out-Null

# Make PayNowNoCFCorp the active tab (mirrors activeTab="2" / tabSelected on sheet3).
$ws3.Activate()
